# Apply weekly report regeneration edits:
# - Update the "Report Generated On" timestamp
# - Zero out the Total Billed Amount and all per-line/day "Pricing" figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"

# Total Billed Amount (Report Summary)
$ws.Range("C8").Value = 0

# Thursday (07/24/2025) block - Pricing column + day total
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0

# Friday (07/25/2025) block - Pricing column + day total
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0

# Saturday (07/26/2025) block - Pricing column + day total
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0

# Sunday (07/27/2025) block - Pricing column + day total
$ws.Range("H49").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("H54").Value = 0
$ws.Range("H55").Value = 0

Write-Output "Applied report regeneration edits (timestamp + zeroed pricing)."
